$wb = $excel.ActiveWorkbook

# Update both the "展览" sheet and the "全部类型" sheet, which carry the
# same data and need the same updated "想去人数" (F column) figures.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1477
    $ws.Range("F4").Value = 61
    $ws.Range("F5").Value = 2224
    $ws.Range("F7").Value = 1366
    $ws.Range("F11").Value = 327
}
